$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.71
$ws.Range("A6").Value = -22.255
$ws.Range("A7").Value = -19.575
$ws.Range("C7").Value = -12.489
$ws.Range("A8").Value = -22.114
$ws.Range("C11").Value = -12.924
$ws.Range("C12").Value = -10.782
$ws.Range("E12").Value = 17.515
$ws.Range("E13").Value = 16.619
$ws.Range("E14").Value = 17.111
$ws.Range("C15").Value = -13.381
$ws.Range("A16").Value = -21.802
$ws.Range("E16").Value = 16.725
$ws.Range("E19").Value = 16.45
$ws.Range("A20").Value = -19.86
$ws.Range("C20").Value = -12.366
$ws.Range("E20").Value = 16.405
$ws.Range("A21").Value = -20.428
$ws.Range("C21").Value = -12.761
$ws.Range("C22").Value = -12.624
$ws.Range("E22").Value = 16.819
$ws.Range("C23").Value = -12.616
$ws.Range("A28").Value = -22.033
$ws.Range("A29").Value = -21.531
$ws.Range("C29").Value = -12.554
$ws.Range("A30").Value = -22.069
$ws.Range("A32").Value = -21.864
$ws.Range("C34").Value = -11.955
$ws.Range("E36").Value = 16.633
$ws.Range("A40").Value = -20.132
$ws.Range("C42").Value = -12.492
$ws.Range("C43").Value = -12.634
$ws.Range("E43").Value = 17.325
$ws.Range("C44").Value = -13.116
$ws.Range("C45").Value = -13.054
$ws.Range("A46").Value = -21.995
$ws.Range("C46").Value = -13.421
$ws.Range("E46").Value = 16.835
$ws.Range("C50").Value = -13.957
$ws.Range("E50").Value = 16.434
$ws.Range("A51").Value = -21.591
$ws.Range("C51").Value = -11.06
$ws.Range("A52").Value = -21.933
$ws.Range("A57").Value = -21.859
$ws.Range("C57").Value = -13.322
$ws.Range("A59").Value = -22.28
$ws.Range("A62").Value = -22.151
$ws.Range("C65").Value = -12.4
$ws.Range("A66").Value = -21.591
$ws.Range("C66").Value = -11.405
$ws.Range("C67").Value = -11.354
$ws.Range("A73").Value = -20.685
$ws.Range("A74").Value = -21.244
$ws.Range("E76").Value = 16.972
$ws.Range("A77").Value = -20.371
$ws.Range("C79").Value = -12.133
$ws.Range("C84").Value = -13.704
$ws.Range("C87").Value = -13.759
$ws.Range("A92").Value = -21.609
$ws.Range("C92").Value = -11.383
$ws.Range("E95").Value = 17.271
$ws.Range("C97").Value = -11.87
$ws.Range("E97").Value = 17.197
$ws.Range("E99").Value = 16.949
$ws.Range("A100").Value = -22.122
